$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Range("H33").Value = 287.6
$ws.Range("I33").Value = 237
$ws.Range("J33").Value = 490
$ws.Range("K33").Value = 237
$ws.Range("L33").Value = 490
$ws.Range("M33").Value = -8
$ws.Range("N33").Value = -948

# Row 80: Cleansing the Wicked Humours | Hallowed Water
$ws.Range("H80").Value = 647.73334
$ws.Range("I80").Value = 326.42856
$ws.Range("J80").Value = 928.875
$ws.Range("K80").Value = 979.28568
$ws.Range("L80").Value = 2786.625
$ws.Range("M80").Value = 18.71432000000004
$ws.Range("N80").Value = -4782.625

# Row 83: Washing Away the Sins (L) | Hallowed Water
$ws.Range("H83").Value = 647.73334
$ws.Range("I83").Value = 326.42856
$ws.Range("J83").Value = 928.875
$ws.Range("K83").Value = 2937.85704
$ws.Range("L83").Value = 8359.875
$ws.Range("M83").Value = 2054.14296
$ws.Range("N83").Value = -18343.875

# Row 100: Asking for a Friend | Beetle Glue
$ws.Range("H100").Value = 2765.875
$ws.Range("I100").Value = 2661.5557
$ws.Range("J100").Value = 2900
$ws.Range("K100").Value = 2661.5557
$ws.Range("L100").Value = 2900
$ws.Range("M100").Value = -2120.5557
$ws.Range("N100").Value = -3982

# Row 132: Fast-forwarding Flora | Growth Formula Lambda
$ws.Range("H132").Value = 894408.8
$ws.Range("I132").Value = 3391.975
$ws.Range("K132").Value = 10175.925
$ws.Range("M132").Value = -7645.924999999999

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3078802.5
$ws.Range("I138").Value = 1919.138
$ws.Range("J138").Value = 5557403
$ws.Range("K138").Value = 5757.414
$ws.Range("L138").Value = 16672209
$ws.Range("M138").Value = -617.4139999999998
$ws.Range("N138").Value = -16682489

# Row 141: Remedy for Reason | Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 5142.625
$ws.Range("I141").Value = 5570.727
$ws.Range("J141").Value = 4200.8
$ws.Range("K141").Value = 16712.181
$ws.Range("L141").Value = 12602.4
$ws.Range("M141").Value = -11532.181
$ws.Range("N141").Value = -22962.4

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks | Mythril Ingot
$ws.Range("H45").Value = 1409.2
$ws.Range("I45").Value = 819.4
$ws.Range("J45").Value = 1999
$ws.Range("K45").Value = 819.4
$ws.Range("L45").Value = 1999
$ws.Range("M45").Value = -442.4
$ws.Range("N45").Value = -2753

# Row 74: As the Bolt Flies | Titanium Nugget
$ws.Range("H74").Value = 5006870.5
$ws.Range("I74").Value = 5766080.5
$ws.Range("J74").Value = 114184.445
$ws.Range("K74").Value = 5766080.5
$ws.Range("L74").Value = 114184.445
$ws.Range("M74").Value = -5765206.5
$ws.Range("N74").Value = -115932.445

# Row 77: Heavy Metal Banned (L) | Titanium Nugget
$ws.Range("H77").Value = 5006870.5
$ws.Range("I77").Value = 5766080.5
$ws.Range("J77").Value = 114184.445
$ws.Range("K77").Value = 28830402.5
$ws.Range("L77").Value = 570922.2250000001
$ws.Range("M77").Value = -28826034.5
$ws.Range("N77").Value = -579658.2250000001

# Row 122: Haste for High Durium | High Durium Nugget
$ws.Range("H122").Value = 3066.25
$ws.Range("I122").Value = 2236.6667
$ws.Range("J122").Value = 5555
$ws.Range("K122").Value = 6710.000100000001
$ws.Range("L122").Value = 16665
$ws.Range("M122").Value = -4260.000100000001
$ws.Range("N122").Value = -21565

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin | Adamantite Nugget
$ws.Range("H86").Value = 10947.615
$ws.Range("I86").Value = 12120.818
$ws.Range("K86").Value = 12120.818
$ws.Range("M86").Value = -10997.818

# Row 89: Piercing Eyes Deserve Piercing Shafts (L) | Adamantite Nugget
$ws.Range("H89").Value = 10947.615
$ws.Range("I89").Value = 12120.818
$ws.Range("K89").Value = 60604.09
$ws.Range("M89").Value = -54988.09

# Row 94: High Steal | High Steel Nugget
$ws.Range("H94").Value = 734.5789
$ws.Range("I94").Value = 732.6429000000001
$ws.Range("K94").Value = 732.6429000000001
$ws.Range("M94").Value = -281.6429000000001

# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Range("H134").Value = 2129.9812
$ws.Range("I134").Value = 2112.9512
$ws.Range("J134").Value = 2188.1667
$ws.Range("K134").Value = 6338.8536
$ws.Range("L134").Value = 6564.500100000001
$ws.Range("M134").Value = -3803.8536
$ws.Range("N134").Value = -11634.5001

$ws = $wb.Worksheets.Item("CRP")
# Row 133: Yimepi's Country Charms | Ginseng Earrings
$ws.Range("H133").Value = 45561.76
$ws.Range("J133").Value = 45561.76
$ws.Range("L133").Value = 45561.76
$ws.Range("N133").Value = -50621.76

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Range("H134").Value = 52667.695
$ws.Range("I134").Value = 6143.1
$ws.Range("J134").Value = 362831.66
$ws.Range("K134").Value = 18429.3
$ws.Range("L134").Value = 1088494.98
$ws.Range("M134").Value = -15894.3
$ws.Range("N134").Value = -1093564.98

# Row 135: The Wing's Wings | Ceiba Wings
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 3: Trout Fishing in Limsa | Grilled Trout
$ws.Range("H3").Value = 8908.725
$ws.Range("I3").Value = 4784.643
$ws.Range("J3").Value = 12757.866
$ws.Range("K3").Value = 14353.929
$ws.Range("L3").Value = 38273.598
$ws.Range("M3").Value = -14241.929
$ws.Range("N3").Value = -38497.598

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Range("H131").Value = 921.56525
$ws.Range("J131").Value = 995.7049
$ws.Range("L131").Value = 2987.1147
$ws.Range("N131").Value = -13067.1147

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence | Ametrine
$ws.Range("H122").Value = 2043
$ws.Range("I122").Value = 2212.5
$ws.Range("K122").Value = 6637.5
$ws.Range("M122").Value = -4187.5

# Row 123: Workplace Workout | Ametrine Ring of Fending
$ws.Range("H123").Value = 39135.668
$ws.Range("J123").Value = 39135.668
$ws.Range("L123").Value = 39135.668
$ws.Range("N123").Value = -44035.668

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Range("H22").Value = 622
$ws.Range("I22").Value = 390.9091
$ws.Range("J22").Value = 1186.8889
$ws.Range("K22").Value = 390.9091
$ws.Range("L22").Value = 1186.8889
$ws.Range("M22").Value = -95.90910000000002
$ws.Range("N22").Value = -1776.8889

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Range("H27").Value = 622
$ws.Range("I27").Value = 390.9091
$ws.Range("J27").Value = 1186.8889
$ws.Range("K27").Value = 390.9091
$ws.Range("L27").Value = 1186.8889
$ws.Range("M27").Value = -283.9091
$ws.Range("N27").Value = -1400.8889

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Range("H136").Value = 51157.547
$ws.Range("I136").Value = 33557.938
$ws.Range("J136").Value = 115689.445
$ws.Range("K136").Value = 100673.814
$ws.Range("L136").Value = 347068.335
$ws.Range("M136").Value = -98123.81400000001
$ws.Range("N136").Value = -352168.335

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke | Rainbow Cloth
$ws.Range("H62").Value = 6194
$ws.Range("I62").Value = 3999.5
$ws.Range("J62").Value = 8388.5
$ws.Range("K62").Value = 3999.5
$ws.Range("L62").Value = 8388.5
$ws.Range("M62").Value = -3375.5
$ws.Range("N62").Value = -9636.5

# Row 65: Desperate for Diversionaries (L) | Rainbow Cloth
$ws.Range("H65").Value = 6194
$ws.Range("I65").Value = 3999.5
$ws.Range("J65").Value = 8388.5
$ws.Range("K65").Value = 19997.5
$ws.Range("L65").Value = 41942.5
$ws.Range("M65").Value = -16877.5
$ws.Range("N65").Value = -48182.5

# Row 113: A Tender Table | Pixie Floss
$ws.Range("H113").Value = 609.0263
$ws.Range("J113").Value = 370.4
$ws.Range("L113").Value = 1111.2
$ws.Range("N113").Value = -5451.2

# Row 122: Heavy Armoire | Dark Hempen Cloth
$ws.Range("H122").Value = 1594.8823
$ws.Range("I122").Value = 890.8
$ws.Range("J122").Value = 2600.7144
$ws.Range("K122").Value = 2672.4
$ws.Range("L122").Value = 7802.1432
$ws.Range("M122").Value = -222.3999999999996
$ws.Range("N122").Value = -12702.1432

# Row 126: A Polished Purchase | Snow Linen
$ws.Range("H126").Value = 1559
$ws.Range("I126").Value = 1501.3334
$ws.Range("J126").Value = 1905
$ws.Range("K126").Value = 4504.0002
$ws.Range("L126").Value = 5715
$ws.Range("M126").Value = -2034.0002
$ws.Range("N126").Value = -10655

# Row 132: Comfy Cabins | Snow Cotton Cloth
$ws.Range("H132").Value = 118712.766
$ws.Range("I132").Value = 111558.22
$ws.Range("J132").Value = 126761.625
$ws.Range("K132").Value = 334674.66
$ws.Range("L132").Value = 380284.875
$ws.Range("M132").Value = -332144.66
$ws.Range("N132").Value = -385344.875

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Range("H136").Value = 41109.367
$ws.Range("I136").Value = 30833.176
$ws.Range("J136").Value = 60519.945
$ws.Range("K136").Value = 92499.52799999999
$ws.Range("L136").Value = 181559.835
$ws.Range("M136").Value = -89949.52799999999
$ws.Range("N136").Value = -186659.835
